# Add a "Label" column (H) to Sheet1:
#  - H1 header "Label", styled like the other header cells (B1:G1)
#  - H2:H11 and H12:H21 filled with 0 for "Control" rows and 1 for "MDD" rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell, matching the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Data rows: two blocks of 10 rows each (2-11 and 12-21),
# first 5 rows of each block are Control (0), last 5 are MDD (1)
$labels = 0,0,0,0,0,1,1,1,1,1

foreach ($block in 0,1) {
    $startRow = 2 + ($block * 10)
    for ($i = 0; $i -lt 10; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, 8).Value = $labels[$i]
    }
}
